$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 15395.218
$ws.Range("I62").Value = 15945
$ws.Range("J62").Value = 3300
$ws.Range("K62").Value = 15945
$ws.Range("L62").Value = 3300
$ws.Range("M62").Value = -15321
$ws.Range("N62").Value = -4548

$ws.Range("H65").Value = 15395.218
$ws.Range("I65").Value = 15945
$ws.Range("J65").Value = 3300
$ws.Range("K65").Value = 79725
$ws.Range("L65").Value = 16500
$ws.Range("M65").Value = -76605
$ws.Range("N65").Value = -22740

$ws.Range("H74").Value = 3310.1
$ws.Range("I74").Value = 3311.8823
$ws.Range("K74").Value = 3311.8823
$ws.Range("M74").Value = -2375.8823

$ws.Range("H77").Value = 3310.1
$ws.Range("I77").Value = 3311.8823
$ws.Range("K77").Value = 16559.4115
$ws.Range("M77").Value = -11879.4115

$ws.Range("H100").Value = 4585
$ws.Range("I100").Value = 3043.5
$ws.Range("J100").Value = 20000
$ws.Range("K100").Value = 3043.5
$ws.Range("L100").Value = 20000
$ws.Range("M100").Value = -2502.5
$ws.Range("N100").Value = -21082

$ws.Range("H112").Value = 1273.8462
$ws.Range("I112").Value = 478
$ws.Range("J112").Value = 1463.3334
$ws.Range("K112").Value = 1434
$ws.Range("L112").Value = 4390.0002
$ws.Range("M112").Value = -326
$ws.Range("N112").Value = -6606.0002

$ws.Range("H121").Value = 1513.909
$ws.Range("J121").Value = 1513.909
$ws.Range("L121").Value = 4541.727000000001
$ws.Range("N121").Value = -8035.727000000001

$ws.Range("H137").Value = 44636.957
$ws.Range("I137").Value = 67702.47
$ws.Range("J137").Value = 1389.125
$ws.Range("K137").Value = 203107.41
$ws.Range("L137").Value = 4167.375
$ws.Range("M137").Value = -200557.41
$ws.Range("N137").Value = -9267.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2437.5
$ws.Range("I61").Value = 1758.3334
$ws.Range("J61").Value = 3116.6667
$ws.Range("K61").Value = 1758.3334
$ws.Range("L61").Value = 3116.6667
$ws.Range("M61").Value = -1546.3334
$ws.Range("N61").Value = -3540.6667

$ws.Range("H74").Value = 38575.89
$ws.Range("I74").Value = 43179.125
$ws.Range("J74").Value = 1750
$ws.Range("K74").Value = 43179.125
$ws.Range("L74").Value = 1750
$ws.Range("M74").Value = -42305.125
$ws.Range("N74").Value = -3498

$ws.Range("H77").Value = 38575.89
$ws.Range("I77").Value = 43179.125
$ws.Range("J77").Value = 1750
$ws.Range("K77").Value = 215895.625
$ws.Range("L77").Value = 8750
$ws.Range("M77").Value = -211527.625
$ws.Range("N77").Value = -17486

$ws.Range("H97").Value = 827.97144
$ws.Range("I97").Value = 547.58826
$ws.Range("J97").Value = 1092.7778
$ws.Range("K97").Value = 547.58826
$ws.Range("L97").Value = 1092.7778
$ws.Range("M97").Value = -51.58825999999999
$ws.Range("N97").Value = -2084.7778

$ws.Range("H110").Value = 2132.625
$ws.Range("I110").Value = 2008.7142
$ws.Range("K110").Value = 2008.7142
$ws.Range("M110").Value = 36.28580000000011

$ws.Range("H136").Value = 2437.5
$ws.Range("I136").Value = 1758.3334
$ws.Range("J136").Value = 3116.6667
$ws.Range("K136").Value = 5275.0002
$ws.Range("L136").Value = 9350.000100000001
$ws.Range("M136").Value = -2725.0002
$ws.Range("N136").Value = -14450.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 306060.12
$ws.Range("I86").Value = 1780.5454
$ws.Range("J86").Value = 584983.0600000001
$ws.Range("K86").Value = 1780.5454
$ws.Range("L86").Value = 584983.0600000001
$ws.Range("M86").Value = -657.5454
$ws.Range("N86").Value = -587229.0600000001

$ws.Range("H89").Value = 306060.12
$ws.Range("I89").Value = 1780.5454
$ws.Range("J89").Value = 584983.0600000001
$ws.Range("K89").Value = 8902.726999999999
$ws.Range("L89").Value = 2924915.3
$ws.Range("M89").Value = -3286.726999999999
$ws.Range("N89").Value = -2936147.3

$ws.Range("H99").Value = 1295.8334
$ws.Range("I99").Value = 1206.25
$ws.Range("J99").Value = 1475
$ws.Range("K99").Value = 1206.25
$ws.Range("L99").Value = 1475
$ws.Range("M99").Value = 291.75
$ws.Range("N99").Value = -4471

$ws.Range("H107").Value = 636
$ws.Range("I107").Value = 583.8
$ws.Range("J107").Value = 810
$ws.Range("K107").Value = 583.8
$ws.Range("L107").Value = 810
$ws.Range("M107").Value = 1336.2
$ws.Range("N107").Value = -4650

$ws.Range("H134").Value = 34041.312
$ws.Range("I134").Value = 2300.16
$ws.Range("J134").Value = 113394.2
$ws.Range("K134").Value = 6900.48
$ws.Range("L134").Value = 340182.6
$ws.Range("M134").Value = -4365.48
$ws.Range("N134").Value = -345252.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 806.3
$ws.Range("I16").Value = 667.5
$ws.Range("J16").Value = 898.8333
$ws.Range("K16").Value = 667.5
$ws.Range("L16").Value = 898.8333
$ws.Range("M16").Value = -380.5
$ws.Range("N16").Value = -1472.8333

$ws.Range("H31").Value = 8628.125
$ws.Range("I31").Value = 9319.4
$ws.Range("J31").Value = 7876.7393
$ws.Range("K31").Value = 9319.4
$ws.Range("L31").Value = 7876.7393
$ws.Range("M31").Value = -9024.4
$ws.Range("N31").Value = -8466.739300000001

$ws.Range("H34").Value = 8628.125
$ws.Range("I34").Value = 9319.4
$ws.Range("J34").Value = 7876.7393
$ws.Range("K34").Value = 9319.4
$ws.Range("L34").Value = 7876.7393
$ws.Range("M34").Value = -9117.4
$ws.Range("N34").Value = -8280.739300000001

$ws.Range("H107").Value = 518.125
$ws.Range("I107").Value = 435.8
$ws.Range("J107").Value = 655.3333
$ws.Range("K107").Value = 435.8
$ws.Range("L107").Value = 655.3333
$ws.Range("M107").Value = 1484.2
$ws.Range("N107").Value = -4495.3333

$ws.Range("H113").Value = 806.3
$ws.Range("I113").Value = 667.5
$ws.Range("J113").Value = 898.8333
$ws.Range("K113").Value = 667.5
$ws.Range("L113").Value = 898.8333
$ws.Range("M113").Value = 1502.5
$ws.Range("N113").Value = -5238.8333

$ws.Range("H132").Value = 43914
$ws.Range("I132").Value = 60593.65
$ws.Range("J132").Value = 3406.2856
$ws.Range("K132").Value = 181780.95
$ws.Range("L132").Value = 10218.8568
$ws.Range("M132").Value = -179250.95
$ws.Range("N132").Value = -15278.8568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 25669070
$ws.Range("I131").Value = 502
$ws.Range("J131").Value = 27196962
$ws.Range("K131").Value = 1506
$ws.Range("L131").Value = 81590886
$ws.Range("M131").Value = 3534
$ws.Range("N131").Value = -81600966

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 810.5217
$ws.Range("I97").Value = 797.61536
$ws.Range("K97").Value = 797.61536
$ws.Range("M97").Value = -301.61536

$ws.Range("H113").Value = 1666.6666
$ws.Range("I113").Value = 1500
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1500
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 670
$ws.Range("N113").Value = -6340

$ws.Range("H132").Value = 55024.895
$ws.Range("I132").Value = 1613.125
$ws.Range("J132").Value = 93869.82000000001
$ws.Range("K132").Value = 4839.375
$ws.Range("L132").Value = 281609.46
$ws.Range("M132").Value = -2309.375
$ws.Range("N132").Value = -286669.46

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 30915.229
$ws.Range("I100").Value = 126763.625
$ws.Range("J100").Value = 2515.7036
$ws.Range("K100").Value = 126763.625
$ws.Range("L100").Value = 2515.7036
$ws.Range("M100").Value = -126222.625
$ws.Range("N100").Value = -3597.7036

$ws.Range("H136").Value = 324719.22
$ws.Range("I136").Value = 626802.9
$ws.Range("J136").Value = 2496.6667
$ws.Range("K136").Value = 1880408.7
$ws.Range("L136").Value = 7490.000100000001
$ws.Range("M136").Value = -1877858.7
$ws.Range("N136").Value = -12590.0001
